$wb = $excel.ActiveWorkbook

# --- Sheet "isa_template": bump template version ---
$wsTemplate = $wb.Worksheets.Item("isa_template")
$wsTemplate.Range("B4").Value = "1.0.1"

# --- Sheet "Tabelle1": annotation table ---
$ws = $wb.Worksheets.Item("Tabelle1")

# Fix up the two header names that previously had an empty MIAPPE id
$ws.Range("AD1").Value = "Term Source REF (MIAPPE:0124)"
$ws.Range("AE1").Value = "Term Accession Number (MIAPPE:0124)"

# Replace the descriptive placeholder row (row 2) with concrete example values
$ws.Range("A2").Value = ""
$ws.Range("B2").Value = "hydroponic plant culture media"
$ws.Range("C2").Value = "EO"
$ws.Range("D2").Value = "http://purl.obolibrary.org/obo/EO_0007067"
$ws.Range("E2").Value = "plant pot"
$ws.Range("F2").Value = "ENVO"
$ws.Range("G2").Value = "http://purl.obolibrary.org/obo/ENVO_03600045"
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = "liter"
$ws.Range("J2").Value = "UO"
$ws.Range("K2").Value = "http://purl.obolibrary.org/obo/UO_0000099"
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = "meter"
$ws.Range("N2").Value = "UO"
$ws.Range("O2").Value = "http://purl.obolibrary.org/obo/UO_0000008"
$ws.Range("P2").Value = "2/container"
$ws.Range("Q2").Value = ""
$ws.Range("R2").Value = ""
$ws.Range("S2").Value = ""
$ws.Range("T2").Value = ""
$ws.Range("U2").Value = ""
$ws.Range("V2").Value = "10/plot"
$ws.Range("W2").Value = ""
$ws.Range("X2").Value = ""
$ws.Range("Y2").Value = "weekly"
$ws.Range("Z2").Value = ""
$ws.Range("AA2").Value = ""
$ws.Range("AB2").Value = "7.7:40-60; 6.5; 4.3:10-20"
$ws.Range("AC2").Value = ""
$ws.Range("AD2").Value = ""
$ws.Range("AE2").Value = ""
$ws.Range("AF2").Value = ""
$ws.Range("AG2").Value = "percent"
$ws.Range("AH2").Value = "UO"
$ws.Range("AI2").Value = "http://purl.obolibrary.org/obo/UO_0000187"
$ws.Range("AJ2").Value = "20"
$ws.Range("AK2").Value = "degree celsius"
$ws.Range("AL2").Value = "UO"
$ws.Range("AM2").Value = "http://purl.obolibrary.org/obo/UO_0000027"
$ws.Range("AN2").Value = "Pa m-2"
$ws.Range("AO2").Value = ""
$ws.Range("AP2").Value = ""
$ws.Range("AQ2").Value = ""
$ws.Range("AR2").Value = "g g-1 dry weight"
$ws.Range("AS2").Value = ""
$ws.Range("AT2").Value = ""
$ws.Range("AU2").Value = ""
$ws.Range("AV2").Value = "percent"
$ws.Range("AW2").Value = ""
$ws.Range("AX2").Value = ""
$ws.Range("AY2").Value = ""
$ws.Range("AZ2").Value = ""

# Drop the old rows 3 and 4 (units/type rows) entirely, shrinking the table to 1 data row
$ws.Rows("3:4").Delete()
